# Rename the header row columns from "_old"/"_new" suffixes to
# "_FV2210"/"_FV2304" suffixes (the two EDI "Formatversion" names being
# diffed), then wrap the data range in a native Excel Table and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10) carry the "_FV2210" (previously "_old") headers,
# columns L-U (12-21) carry the "_FV2304" (previously "_new") headers.
# Column K ("diff") is untouched.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $leftCol = $i + 1
    $rightCol = $i + 12
    $ws.Cells.Item(1, $leftCol).Value = $baseNames[$i] + "_FV2210"
    $ws.Cells.Item(1, $rightCol).Value = $baseNames[$i] + "_FV2304"
}

# Turn the used range into a real Excel Table (ListObject) so the header
# row doubles as filter buttons, picking up the renamed headers above as
# its column names.
$dataRange = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row (split after row 1, keep column A as the left edge).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
